$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Depth"
$ws.Range("C1").Value = "DRM"
$ws.Range("D1").Value = "Month"
$ws.Range("E1").Value = "Depth:DRM"
$ws.Range("F1").Value = "Depth:Month"
$ws.Range("G1").Value = "DRM:Month"
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = -0.3395471444682089
$ws.Range("D2").Value = "+"
$ws.Range("F2").ClearContents()
$ws.Range("G2").Value = "+"
$ws.Range("H2").Value = 0.330731964963923
$ws.Range("K2").Value = 63.88230265466321
$ws.Range("M2").Value = 0.4893838242007741
$ws.Range("N2").Value = 0.2614973406498461
$ws.Range("B3").Value = 0.07911341566539438
$ws.Range("C3").Value = -0.3026235179555763
$ws.Range("D3").Value = "+"
$ws.Range("F3").ClearContents()
$ws.Range("G3").Value = "+"
$ws.Range("H3").Value = 0.3425687334411727
$ws.Range("L3").Value = 2.419681507154564
$ws.Range("M3").Value = 0.1459561661204095
$ws.Range("N3").Value = 0.2486499810756259
$ws.Range("B4").Value = -0.1438656651210444
$ws.Range("C4").Value = -0.4066917879176565
$ws.Range("D4").Value = "+"
$ws.Range("E4").ClearContents()
$ws.Range("G4").Value = "+"
$ws.Range("J4").Value = -24.00158869585068
$ws.Range("K4").Value = 66.48317739170136
$ws.Range("L4").Value = 2.600874737038154
$ws.Range("M4").Value = 0.1333143308470349
$ws.Range("B5").Value = 0.08694104806002551
$ws.Range("C5").Value = -0.2328601803217152
$ws.Range("D5").Value = "+"
$ws.Range("E5").Value = 0.1486646389153216
$ws.Range("F5").ClearContents()
$ws.Range("G5").Value = "+"
$ws.Range("H5").Value = 0.3635163605442524
$ws.Range("K5").Value = 68.48262519061171
$ws.Range("M5").Value = 0.04905714436676108
$ws.Range("N5").Value = 0.2456490199042991
$ws.Range("J6").Value = -32.45595592166235
$ws.Range("K6").Value = 69.31191184332471
$ws.Range("L6").Value = 5.429609188661502
$ws.Range("M6").Value = 0.03240596477422415
$ws.Range("B7").Value = 0.1771497258443417
$ws.Range("C7").ClearContents()
$ws.Range("H7").Value = 0.06815431984537323
$ws.Range("L7").Value = 5.527789506653434
$ws.Range("M7").Value = 0.03085356629572889
$ws.Range("N7").Value = 0.03809478177586911
$ws.Range("B8").Value = -0.1178129235351062
$ws.Range("C8").Value = -0.362097157316762
$ws.Range("D8").Value = "+"
$ws.Range("E8").Value = 0.07293885038244627
$ws.Range("G8").Value = "+"
$ws.Range("H8").Value = 0.405398637733311
$ws.Range("L8").Value = 5.874094635890202
$ws.Range("M8").Value = 0.02594813609395962
$ws.Range("N8").Value = 0.2681829387486905
$ws.Range("C9").Value = 0.1170974647800254
$ws.Range("D9").ClearContents()
$ws.Range("J9").Value = -31.95714100363683
$ws.Range("K9").Value = 70.74186821417022
$ws.Range("L9").Value = 6.859565559507018
$ws.Range("M9").Value = 0.01585308816734456
$ws.Range("B10").Value = 0.1833539723030487
$ws.Range("C10").Value = 0.1261174181684735
$ws.Range("D10").ClearContents()
$ws.Range("H10").Value = 0.1026139224622499
$ws.Range("J10").Value = -30.66951577525997
$ws.Range("K10").Value = 70.76760297909136
$ws.Range("M10").Value = 0.01565040720193887
$ws.Range("N10").Value = 0.04278818395973316
$ws.Range("B11").Value = 0.1893523149359384
$ws.Range("C11").ClearContents()
$ws.Range("D11").Value = "+"
$ws.Range("H11").Value = 0.0943372137216959
$ws.Range("L11").Value = 7.188268604872277
$ws.Range("M11").Value = 0.01345045313370756
$ws.Range("N11").Value = 0.03395969463647563
$ws.Range("B12").ClearContents()
$ws.Range("D12").Value = "+"
$ws.Range("H12").Value = 0.01742017036698895
$ws.Range("L12").Value = 7.277263724407021
$ws.Range("M12").Value = 0.0128650616262989
$ws.Range("N12").Value = -0.01427595316955976
$ws.Range("B13").Value = 0.197265482883538
$ws.Range("C13").Value = 0.2058344846318508
$ws.Range("D13").ClearContents()
$ws.Range("E13").Value = 0.2197762160638421
$ws.Range("G13").ClearContents()
$ws.Range("H13").Value = 0.1600601718682541
$ws.Range("L13").Value = 7.495805655611093
$ws.Range("M13").Value = 0.01153336651424878
$ws.Range("N13").Value = 0.07316984482014244
$ws.Range("B14").Value = 0.03338127796077408
$ws.Range("C14").ClearContents()
$ws.Range("D14").Value = "+"
$ws.Range("E14").ClearContents()
$ws.Range("F14").Value = "+"
$ws.Range("J14").Value = -30.10629212519849
$ws.Range("K14").Value = 72.4348064726192
$ws.Range("L14").Value = 8.552503817955994
$ws.Range("M14").Value = 0.006799813526122207
$ws.Range("B15").Value = 0.1946976382412989
$ws.Range("C15").Value = 0.1212986907912633
$ws.Range("D15").Value = "+"
$ws.Range("H15").Value = 0.1261458393100859
$ws.Range("L15").Value = 8.802052623061314
$ws.Range("M15").Value = 0.006002168296075466
$ws.Range("N15").Value = 0.03574713303181898
$ws.Range("B16").ClearContents()
$ws.Range("C16").Value = 0.1128011916503565
$ws.Range("D16").Value = "+"
$ws.Range("H16").Value = 0.04498124441349627
$ws.Range("L16").Value = 8.939378739123114
$ws.Range("M16").Value = 0.005603871670534263
$ws.Range("N16").Value = -0.01868667262560408
$ws.Range("B17").Value = 0.2001267315573115
$ws.Range("C17").Value = 0.1952912399265134
$ws.Range("D17").Value = "+"
$ws.Range("E17").Value = 0.1958347189466294
$ws.Range("G17").ClearContents()
$ws.Range("H17").Value = 0.1628815042493786
$ws.Range("L17").Value = 10.39332008397328
$ws.Range("M17").Value = 0.002708749410304898
$ws.Range("N17").Value = 0.04329314771357562
$ws.Range("B18").Value = 0.07102380370289117
$ws.Range("C18").Value = 0.08637049434884325
$ws.Range("D18").Value = "+"
$ws.Range("E18").ClearContents()
$ws.Range("F18").Value = "+"
$ws.Range("H18").Value = 0.1470496042877665
$ws.Range("L18").Value = 11.01160018400358
$ws.Range("M18").Value = 0.001988433228260608
$ws.Range("N18").Value = 0.02519954775744737
$ws.Range("B19").Value = 0.1187194610495932
$ws.Range("C19").Value = 0.1612557501989792
$ws.Range("D19").Value = "+"
$ws.Range("E19").Value = 0.1659854395707197
$ws.Range("F19").Value = "+"
$ws.Range("G19").ClearContents()
$ws.Range("H19").Value = 0.170902086865355
$ws.Range("J19").Value = -29.36357507700962
$ws.Range("K19").Value = 77.20715015401925
$ws.Range("L19").Value = 13.32484749935604
$ws.Range("M19").Value = 0.0006254545262718938
$ws.Range("N19").Value = 0.0173654362848652
